# S-01015: se cargan las horas insumidas para parsear el archivo HF.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Horas insumidas")

# Insert a new row before row 44 (shifts existing rows 44-64 down to 45-65,
# and Excel auto-adjusts dependent formulas/ranges).
$ws.Rows("44:44").Insert()

# The freshly inserted row inherits formatting from the row above; reset it
# to plain (unstyled) cells, matching the rest of the data rows.
$ws.Range("B44:F44").ClearFormats()
$ws.Range("B44").NumberFormat = $ws.Range("B45").NumberFormat

# Fill in the new row with the "Proceso archivo HF" entry for story S-01015.
$ws.Range("B44").Value = 40450
$ws.Range("C44").Value = "Duilio"
$ws.Range("E44").Value = "S-01015"
$ws.Range("D44").Value = "Proceso archivo HF"
$ws.Range("F44").Value = 4

# Restore the view state shown in the edited workbook (best effort; some
# hosts may not expose window scroll position).
try {
    $ws.Application.ActiveWindow.ScrollRow = 31
} catch {
}
$ws.Range("E46").Select()
